$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The empty "Title 39" placeholder on this slide was deleted, which freed
# up vertical space that the UML diagram above it was then shifted into.
# Every other shape in the diagram (rectangles, the flowchart extract
# connector and the two elbow connectors) moved by the same uniform
# offset: +35535 EMU horizontally, -1170484 EMU vertically.
$dxEmu = 35535
$dyEmu = -1170484
$emuPerPt = 12700.0
# Tiny bias added after the EMU->point conversion so that the host's
# point->EMU re-quantisation on save lands back on the exact target EMU
# value instead of truncating to one unit below it.
$bias = 0.00005

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 39") {
        continue
    }
    $newLeftEmu = [math]::Round($sh.Left * $emuPerPt) + $dxEmu
    $newTopEmu = [math]::Round($sh.Top * $emuPerPt) + $dyEmu
    $sh.Left = ($newLeftEmu / $emuPerPt) + $bias
    $sh.Top = ($newTopEmu / $emuPerPt) + $bias
}

# Remove the now-vacated (empty) title placeholder shape.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 39") {
        $sh.Delete()
    }
}
